$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 155.8
$ws.Range("J9").Value = 58.666668
$ws.Range("L9").Value = 58.666668
$ws.Range("N9").Value = -396.666668

$ws.Range("H53").Value = 274.2857
$ws.Range("I53").Value = 153.33333
$ws.Range("K53").Value = 153.33333
$ws.Range("M53").Value = 483.66667

$ws.Range("H55").Value = 446.9375
$ws.Range("I55").Value = 127.15385
$ws.Range("J55").Value = 1832.6666
$ws.Range("K55").Value = 127.15385
$ws.Range("L55").Value = 1832.6666
$ws.Range("M55").Value = 86.84614999999999
$ws.Range("N55").Value = -2260.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2365.4285
$ws.Range("I74").Value = 2258.7
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 2258.7
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -1384.7
$ws.Range("N74").Value = -6248

$ws.Range("H77").Value = 2365.4285
$ws.Range("I77").Value = 2258.7
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 11293.5
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -6925.5
$ws.Range("N77").Value = -31236

$ws.Range("H88").Value = 3006.6
$ws.Range("I88").Value = 3503.3333
$ws.Range("J88").Value = 2261.5
$ws.Range("K88").Value = 3503.3333
$ws.Range("L88").Value = 2261.5
$ws.Range("M88").Value = -3097.3333
$ws.Range("N88").Value = -3073.5

$ws.Range("H91").Value = 3006.6
$ws.Range("I91").Value = 3503.3333
$ws.Range("J91").Value = 2261.5
$ws.Range("K91").Value = 3503.3333
$ws.Range("L91").Value = 2261.5
$ws.Range("M91").Value = -2099.3333
$ws.Range("N91").Value = -5069.5

$ws.Range("H97").Value = 658.4666999999999
$ws.Range("I97").Value = 658.4666999999999
$ws.Range("K97").Value = 658.4666999999999
$ws.Range("M97").Value = -162.4666999999999

$ws.Range("H130").Value = 20426.5
$ws.Range("J130").Value = 20426.5
$ws.Range("L130").Value = 20426.5
$ws.Range("N130").Value = -30466.5

$ws.Range("H132").Value = 3129.318
$ws.Range("I132").Value = 3078.3
$ws.Range("K132").Value = 9234.900000000001
$ws.Range("M132").Value = -6704.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7771.1816
$ws.Range("I86").Value = 4997.5
$ws.Range("J86").Value = 9356.143
$ws.Range("K86").Value = 4997.5
$ws.Range("L86").Value = 9356.143
$ws.Range("M86").Value = -3874.5
$ws.Range("N86").Value = -11602.143

$ws.Range("H89").Value = 7771.1816
$ws.Range("I89").Value = 4997.5
$ws.Range("J89").Value = 9356.143
$ws.Range("K89").Value = 24987.5
$ws.Range("L89").Value = 46780.715
$ws.Range("M89").Value = -19371.5
$ws.Range("N89").Value = -58012.715

$ws.Range("H94").Value = 1359.3043
$ws.Range("I94").Value = 1402.95
$ws.Range("J94").Value = 1068.3334
$ws.Range("K94").Value = 1402.95
$ws.Range("L94").Value = 1068.3334
$ws.Range("M94").Value = -951.95
$ws.Range("N94").Value = -1970.3334

$ws.Range("H134").Value = 2499
$ws.Range("I134").Value = 2499
$ws.Range("K134").Value = 7497
$ws.Range("M134").Value = -4962

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3626.6206
$ws.Range("I7").Value = 5783.3335
$ws.Range("J7").Value = 97.454544
$ws.Range("K7").Value = 5783.3335
$ws.Range("L7").Value = 97.454544
$ws.Range("M7").Value = -5670.3335
$ws.Range("N7").Value = -323.454544

$ws.Range("H58").Value = 3523.4614
$ws.Range("I58").Value = 3081.2
$ws.Range("K58").Value = 3081.2
$ws.Range("M58").Value = -2878.2

$ws.Range("H99").Value = 2397.5
$ws.Range("J99").Value = 2397.5
$ws.Range("L99").Value = 2397.5
$ws.Range("N99").Value = -5393.5

$ws.Range("H103").Value = 18643
$ws.Range("I103").Value = 18643
$ws.Range("K103").Value = 18643
$ws.Range("M103").Value = -17471

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws.Range("H126").Value = 2397.5
$ws.Range("J126").Value = 2397.5
$ws.Range("L126").Value = 7192.5
$ws.Range("N126").Value = -12132.5

$ws.Range("H132").Value = 2959.0625
$ws.Range("I132").Value = 2488.077
$ws.Range("K132").Value = 7464.231000000001
$ws.Range("M132").Value = -4934.231000000001

$ws.Range("H134").Value = 2006
$ws.Range("I134").Value = 2006
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6018
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3483
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 3523.4614
$ws.Range("I136").Value = 3081.2
$ws.Range("K136").Value = 9243.599999999999
$ws.Range("M136").Value = -6693.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 129436.25
$ws.Range("J4").Value = 5453.636
$ws.Range("L4").Value = 16360.908
$ws.Range("N4").Value = -16584.908

$ws.Range("H37").Value = 107947.5
$ws.Range("J37").Value = 107947.5
$ws.Range("L37").Value = 323842.5
$ws.Range("N37").Value = -324066.5

$ws.Range("H68").Value = 2057.7778
$ws.Range("I68").Value = 1515
$ws.Range("J68").Value = 2492
$ws.Range("K68").Value = 4545
$ws.Range("L68").Value = 7476
$ws.Range("M68").Value = -3734
$ws.Range("N68").Value = -9098

$ws.Range("H71").Value = 2057.7778
$ws.Range("I71").Value = 1515
$ws.Range("J71").Value = 2492
$ws.Range("K71").Value = 13635
$ws.Range("L71").Value = 22428
$ws.Range("M71").Value = -9579
$ws.Range("N71").Value = -30540

$ws.Range("H104").Value = 9158.25
$ws.Range("I104").Value = 5000
$ws.Range("J104").Value = 9989.9
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 29969.7
$ws.Range("M104").Value = -12379
$ws.Range("N104").Value = -35211.7

$ws.Range("H109").Value = 41202.832
$ws.Range("J109").Value = 1630
$ws.Range("L109").Value = 4890
$ws.Range("N109").Value = -6970

$ws.Range("H114").Value = 976.1111
$ws.Range("I114").Value = 1040.7142
$ws.Range("J114").Value = 750
$ws.Range("K114").Value = 3122.1426
$ws.Range("L114").Value = 2250
$ws.Range("M114").Value = 131.8574000000003
$ws.Range("N114").Value = -8758

$ws.Range("H117").Value = 2899.4
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2899.4
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 8698.200000000001
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -15582.2

$ws.Range("H121").Value = 531.3333
$ws.Range("J121").Value = 955.3333
$ws.Range("L121").Value = 2865.9999
$ws.Range("N121").Value = -5485.9999

$ws.Range("H129").Value = 2232.8
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2232.8
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 6698.400000000001
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -16698.4

$ws.Range("H131").Value = 1541.5834
$ws.Range("I131").Value = 873.5
$ws.Range("K131").Value = 2620.5
$ws.Range("M131").Value = 2419.5

$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 3000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -3900

$ws.Range("H139").Value = 3726.7778
$ws.Range("I139").Value = 3191.5715
$ws.Range("K139").Value = 9574.7145
$ws.Range("M139").Value = -4434.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H80").Value = 13332.333
$ws.Range("J80").Value = 14999.5
$ws.Range("L80").Value = 14999.5
$ws.Range("N80").Value = -16995.5

$ws.Range("H83").Value = 13332.333
$ws.Range("J83").Value = 14999.5
$ws.Range("L83").Value = 74997.5
$ws.Range("N83").Value = -84981.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 8000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 8000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H125").Value = 190357.5
$ws.Range("J125").Value = 190357.5
$ws.Range("L125").Value = 190357.5
$ws.Range("N125").Value = -200197.5

$ws.Range("H132").Value = 1984.5
$ws.Range("I132").Value = 1984.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5953.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3423.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 23726.076
$ws.Range("I2").Value = 23888.334
$ws.Range("K2").Value = 23888.334
$ws.Range("M2").Value = -23776.334

$ws.Range("H126").Value = 5307.0835
$ws.Range("I126").Value = 3428.923
$ws.Range("J126").Value = 7526.727
$ws.Range("K126").Value = 10286.769
$ws.Range("L126").Value = 22580.181
$ws.Range("M126").Value = -7816.769
$ws.Range("N126").Value = -27520.181

$ws.Range("H132").Value = 2259.7144
$ws.Range("I132").Value = 2259.7144
$ws.Range("K132").Value = 6779.1432
$ws.Range("M132").Value = -4249.1432

$ws.Range("H136").Value = 4041.3914
$ws.Range("I136").Value = 2998.0833
$ws.Range("K136").Value = 8994.249899999999
$ws.Range("M136").Value = -6444.249899999999
